$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the "Desarquivamentos Pendentes" sheet entirely
$wsDesarquivamentos = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$wsDesarquivamentos.Delete()

# Deleting the (now) last sheet shifts Excel's active tab; restore the
# original active sheet so tabSelected stays on "PAINEIS DARQ".
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
